$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range('D2').Value = '29.339.17'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '1.876.03'
$ws.Range('E3').Value = '  +0.94%  '
Set-TextValue 'D4' '0.9997'
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '0.7117'
$ws.Range('E5').Value = '  -0.29%  '
Set-TextValue 'D6' '242.47'
$ws.Range('E6').Value = '  +0.92%  '
Set-TextValue 'D7' '1.0000'
$ws.Range('E7').Value = '  -0.07%  '
Set-TextValue 'D8' '0.3113'
$ws.Range('E8').Value = '  +1.25%  '
Set-TextValue 'D9' '0.07756'
$ws.Range('E9').Value = '  +0.32%  '
Set-TextValue 'D10' '25.07'
$ws.Range('E10').Value = '  +0.59%  '
Set-TextValue 'D11' '0.08472'
$ws.Range('E11').Value = '  +2.64%  '
$ws.Range('D12').Value = '1.923.64'
$ws.Range('E12').Value = '  +3.46%  '
Set-TextValue 'D13' '5.208'
$ws.Range('E13').Value = '  -0.19%  '
Set-TextValue 'D14' '0.7111'
$ws.Range('E14').Value = '  -0.42%  '
Set-TextValue 'D15' '91.39'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '29.336.67'
$ws.Range('E16').Value = '  +0.35%  '
Set-TextValue 'D17' '0.000008287'
$ws.Range('E17').Value = '  +6.20%  '
Set-TextValue 'D18' '6.002'
$ws.Range('E18').Value = '  +2.36%  '
Set-TextValue 'D19' '242.73'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D20' '13.22'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.124.49'
$ws.Range('E21').Value = '  +0.53%  '
Set-TextValue 'D22' '0.9996'
$ws.Range('E22').Value = '  -0.15%  '
Set-TextValue 'D23' '7.835'
$ws.Range('E23').Value = '  -1.17%  '
Set-TextValue 'D24' '0.9998'
$ws.Range('E24').Value = '  -0.10%  '
Set-TextValue 'D25' '0.1612'
$ws.Range('E25').Value = '  +2.16%  '
Set-TextValue 'D26' '162.79'
$ws.Range('E26').Value = '  +0.15%  '
Set-TextValue 'D27' '9.024'
$ws.Range('E27').Value = '  +1.45%  '
Set-TextValue 'D28' '18.49'
$ws.Range('E28').Value = '  +1.39%  '
Set-TextValue 'D29' '1.515'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('E30').Value = '  +0.94%  '
Set-TextValue 'D31' '4.350'
$ws.Range('E31').Value = '  +5.81%  '
Set-TextValue 'D32' '1.275'
$ws.Range('E32').Value = '  -3.39%  '
Set-TextValue 'D33' '0.05267'
$ws.Range('E33').Value = '  +1.62%  '
Set-TextValue 'D34' '1.931'
$ws.Range('E34').Value = '  +1.36%  '
Set-TextValue 'D35' '1.175'
$ws.Range('E35').Value = '  +0.22%  '
Set-TextValue 'D36' '0.7443'
$ws.Range('E36').Value = '  +2.30%  '
$ws.Range('E37').Value = '  +0.09%  '
Set-TextValue 'D38' '0.01869'
$ws.Range('E38').Value = '  +1.31%  '
Set-TextValue 'D39' '2.718'
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('D40').Value = '1.171.75'
$ws.Range('E40').Value = '  +1.75%  '
Set-TextValue 'D41' '6.380'
$ws.Range('E41').Value = '  +4.68%  '
Set-TextValue 'D42' '73.23'
$ws.Range('E42').Value = '  +1.60%  '
Set-TextValue 'D43' '0.8880'
$ws.Range('E43').Value = '  -1.48%  '
Set-TextValue 'D44' '106.85'
$ws.Range('E44').Value = '  +4.99%  '
Set-TextValue 'D45' '0.9997'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '2.021.13'
$ws.Range('E46').Value = '  +0.31%  '
Set-TextValue 'D47' '1.817'
$ws.Range('E47').Value = '  +3.11%  '
Set-TextValue 'D48' '0.5203'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '9.386'
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D50' '0.4302'
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D51' '7.059'
$ws.Range('E51').Value = '  +0.57%  '
